# Updated cryptos list on Sat Sep 23 17:50:24 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) are free-text values (e.g. thousand-separated
# strings like "26.707.05") that must stay text. Forcing NumberFormat to
# "@" before assignment keeps Excel from reinterpreting them as numbers
# (which would silently drop meaningful trailing zeros, e.g. "209.30").
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.707.05"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.597.75"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.27"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.51"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.822.22"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.573.92"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.32"
$ws.Range("E17").Value = "  +5.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.674.06"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "209.30"
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.15"
$ws.Range("E21").Value = "  +5.08%  "
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.10"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.31"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("E30").Value = "  +3.12%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.285.60"
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("E35").Value = "  -6.34%  "
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("E39").Value = "  +16.22%  "
$ws.Range("E40").Value = "  -1.13%  "
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.785"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.21"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.734.84"
$ws.Range("E46").Value = "  +1.90%  "
$ws.Range("E47").Value = "  -2.17%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("E51").Value = "  -1.39%  "
